$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 124, shifting the existing
# rows 124-128 down to 126-130 (carrying their formatting/values with them).
$ws.Rows("124:125").Insert()

# New row 124: weekly observation for "Especial" quality.
$ws.Cells.Item(124, 1).Value2 = 7
$ws.Cells.Item(124, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(124, 3).Value2 = "Ñuble"
$ws.Cells.Item(124, 4).Value2 = 44461
$ws.Cells.Item(124, 5).Value2 = 16
$ws.Cells.Item(124, 6).Value2 = "Fruta"
$ws.Cells.Item(124, 7).Value2 = 100104
$ws.Cells.Item(124, 8).Value2 = "Frutos de pepita"
$ws.Cells.Item(124, 9).Value2 = 100104005
$ws.Cells.Item(124, 10).Value2 = "Pera"
$ws.Cells.Item(124, 11).Value2 = "Packham's Triumph"
$ws.Cells.Item(124, 12).Value2 = "Especial"
$ws.Cells.Item(124, 13).Value2 = 60
$ws.Cells.Item(124, 14).Value2 = 11000
$ws.Cells.Item(124, 15).Value2 = 11000
$ws.Cells.Item(124, 16).Value2 = 11000
$ws.Cells.Item(124, 17).Value2 = "$/caja 16 kilos empedrada"
$ws.Cells.Item(124, 18).Value2 = "Provincia de Curicó"
$ws.Cells.Item(124, 19).Value2 = 688
$ws.Cells.Item(124, 20).Value2 = 16

# New row 125: weekly observation for "Primera" quality.
$ws.Cells.Item(125, 1).Value2 = 7
$ws.Cells.Item(125, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(125, 3).Value2 = "Ñuble"
$ws.Cells.Item(125, 4).Value2 = 44461
$ws.Cells.Item(125, 5).Value2 = 16
$ws.Cells.Item(125, 6).Value2 = "Fruta"
$ws.Cells.Item(125, 7).Value2 = 100104
$ws.Cells.Item(125, 8).Value2 = "Frutos de pepita"
$ws.Cells.Item(125, 9).Value2 = 100104005
$ws.Cells.Item(125, 10).Value2 = "Pera"
$ws.Cells.Item(125, 11).Value2 = "Packham's Triumph"
$ws.Cells.Item(125, 12).Value2 = "Primera"
$ws.Cells.Item(125, 13).Value2 = 60
$ws.Cells.Item(125, 14).Value2 = 9000
$ws.Cells.Item(125, 15).Value2 = 10000
$ws.Cells.Item(125, 16).Value2 = 9500
$ws.Cells.Item(125, 17).Value2 = "$/caja 16 kilos empedrada"
$ws.Cells.Item(125, 18).Value2 = "Provincia de Curicó"
$ws.Cells.Item(125, 19).Value2 = 594
$ws.Cells.Item(125, 20).Value2 = 16
